{"js": "// 1) {responsibles} -> {responsible} in the \"Verantwortlich\" row.\nconst respResults = context.document.body.search(\"{responsibles}\", { matchCase: true });\nrespResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < respResults.items.length; i++) {\n  respResults.items[i].insertText(\"{responsible}\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Fill the empty \"Mitwirkend\" (participants) cell with {participants}.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (let t = 0; t < tables.items.length; t++) {\n  const table = tables.items[t];\n  table.load(\"values\");\n}\nawait context.sync();\n\nfor (let t = 0; t < tables.items.length; t++) {\n  const table = tables.items[t];\n  const values = table.values;\n  for (let r = 0; r < values.length; r++) {\n    if (values[r][0] && values[r][0].trim() === \"Mitwirkend\") {\n      const cell = table.getCellOrNullObject(r, 1);\n      await context.sync();\n      cell.body.clear();\n      cell.body.insertText(\"{participants}\", Word.InsertLocation.start);\n      await context.sync();\n    }\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) {responsibles} -> {responsible} in the \"Verantwortlich\" row.\n$find = $d.Content.Find\n$find.Text = \"{responsibles}\"\n$find.Replacement.Text = \"{responsible}\"\n$wdFindContinue = 1\n$wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n\n# 2) Fill the empty \"Mitwirkend\" (participants) cell with {participants}.\nfor ($i = 1; $i -le $d.Tables.Count; $i++) {\n    $table = $d.Tables.Item($i)\n    for ($r = 1; $r -le $table.Rows.Count; $r++) {\n        try {\n            $label = $table.Cell($r, 1).Range.Text\n            $label = $label.TrimEnd([char]13, [char]7).Trim()\n            if ($label -eq \"Mitwirkend\") {\n                $table.Cell($r, 2).Range.Text = \"{participants}\"\n            }\n        } catch {\n            continue\n        }\n    }\n}\n"}
